$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after "Operadores_búsqueda" and rename it
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Ejemplos_Avanzados"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 80.85546875
$ws2.Columns.Item(2).ColumnWidth = 53.42578125
$ws2.Columns.Item(3).ColumnWidth = 52.5703125

# Row 1
$ws2.Range("A1").Value = "¿Cuál sería la expresión correspondiente a la búsqueda de ficheros xls que contengan las palabras user y password dentro de un servidor ftp?"
$ws2.Range("A1").WrapText = $true
$ws2.Range("B1").Value = 'ext:xls inurl:"ftp://" user password'

# Row 2
$ws2.Range("A2").Value = '¿Cómo identificarías páginas gubernamentales que han sido hackeadas? La frase exacta debe ser:"hacked by". '
$ws2.Range("A2").WrapText = $true
$ws2.Range("B2").Value = 'site:gob.* "hacked by"'

# Row 3
$ws2.Range("A3").Value = "¿Cómo identificarías algún sistema SCADA?"
$ws2.Range("B3").Value = 'inurl:"Portal/Portal.mwsl"'
$ws2.Range("B3").WrapText = $true
$ws2.Range("B3").Font.Color = 0x413B37
$ws2.Range("C3").Value = "https://www.flu-project.com/2016/05/google-hacking-de-sistemas-scada-de.html"
$ws2.Range("I3").Value = "https://www.hackers-arise.com/post/2016/07/05/scada-hacking-finding-vulnerable-scada-systems-using-google-hacking"

# Row 4
$ws2.Range("A4").Value = "¿Cómo identificarías en sitios gubernamentales algún fichero ofimático con la marca de confidencial?"
$ws2.Range("A4").WrapText = $true
$ws2.Range("B4").Value = "allintitle:confidential filetype:doc site:gob"
$ws2.Range("C4").Value = [char]0x201C + "robots.txt" + [char]0x201D + " " + [char]0x201C + "disallow:" + [char]0x201D + " filetype:txt site:gob"

# Row 5
$ws2.Range("A5").Value = "¿Podrías identificar algún tipo de vulnerabilidad en alguna página gubernamental peruana?"
$ws2.Range("A5").WrapText = $true

# Row 6
$ws2.Range("A6").Value = "Usando los ejemplos anteriores:"
$ws2.Range("A6").WrapText = $true
$ws2.Range("B6").Value = [char]0x201C + "robots.txt" + [char]0x201D + " " + [char]0x201C + "disallow:" + [char]0x201D + " filetype:txt site:gob.pe"
$ws2.Range("C6").Value = "Se pueden ficheros que no quieren que sean indexados."

# Row 7
$ws2.Range("B7").Value = "allintitle:confidential filetype:pdf site:gob.pe"
$ws2.Range("C7").Value = "Algunos documentos confidenciales"

# Row heights to match target (51.75 for rows 1-5)
$ws2.Rows.Item(1).RowHeight = 51.75
$ws2.Rows.Item(2).RowHeight = 51.75
$ws2.Rows.Item(3).RowHeight = 51.75
$ws2.Rows.Item(4).RowHeight = 51.75
$ws2.Rows.Item(5).RowHeight = 51.75

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Sheet1 selection tweak
[void]$ws1.Range("A9").Select()

[void]$ws2.Select()
[void]$ws2.Range("C8").Select()

Write-Host "done"
